# Trade #23 closed at 2026-02-17 20:07:14 - unknown UNKNOWN +0.000%
#
# Updates the rolled-up Summary / Strategy Status metrics and appends the
# new closed trade (row 24) to both the "All Trades" and "MarketMaking"
# trade logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet - refresh aggregate stats
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1399.83   # Current Capital
$summary.Range("B4").Value = -0.18     # Total P&L $
$summary.Range("B5").Value = -0.16     # Total P&L %
$summary.Range("B6").Value = 23        # Total Trades
$summary.Range("B7").Value = 12        # Winning Trades
$summary.Range("B9").Value = 52.17     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 99.83      # Capital
$status.Range("D5").Value = 23         # Trades
$status.Range("E5").Value = -0.18      # P&L $
$status.Range("F5").Value = -0.17      # P&L %
$status.Range("G5").Value = 52.17      # Win Rate %

# ---------------------------------------------------------------------
# Append the new trade row (row 24) to a trade-log style sheet.
#
# Date/time columns (B, C) hold plain text like "2026-02-17" /
# "20:07:07" in this workbook (not real date/time values), so the
# destination cells are pre-formatted as Text to stop Excel's
# autoconvert-to-date-serial behaviour; the format is reset back to
# Normal immediately after so no stray number format lingers on the
# cell once the literal text is stored.
# ---------------------------------------------------------------------
function Add-TradeRow($sheet) {
    $dateTimeCells = $sheet.Range("B24:C24")
    $dateTimeCells.NumberFormat = "@"

    $sheet.Cells.Item(24, 1).Value = 23
    $sheet.Cells.Item(24, 2).Value = "2026-02-17"
    $sheet.Cells.Item(24, 3).Value = "20:07:07"
    $sheet.Cells.Item(24, 4).Value = "MarketMaking"
    $sheet.Cells.Item(24, 5).Value = "UP"
    $sheet.Cells.Item(24, 6).Value = 0.9
    $sheet.Cells.Item(24, 7).Value = 0.92
    $sheet.Cells.Item(24, 8).Value = "CLOSED"
    $sheet.Cells.Item(24, 9).Value = 2.2222
    $sheet.Cells.Item(24, 10).Value = 0.02
    $sheet.Cells.Item(24, 11).Value = 99.83
    $sheet.Cells.Item(24, 12).Value = 0
    $sheet.Cells.Item(24, 13).Value = 0
    $sheet.Cells.Item(24, 14).Value = 0.6
    $sheet.Cells.Item(24, 15).Value = "Normal spread capture: 19600 bps"
    $sheet.Cells.Item(24, 16).Value = "early_exit"
    $sheet.Cells.Item(24, 17).Value = 0.14

    $dateTimeCells.Style = "Normal"
}

Add-TradeRow $wb.Worksheets.Item("All Trades")
Add-TradeRow $wb.Worksheets.Item("MarketMaking")
